$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.210.32"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "3.793.65"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.18"
$ws.Range("E5").Value = "  +3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.82"
$ws.Range("E6").Value = "  -3.39%  "
$ws.Range("D7").Value = "3.790.85"
$ws.Range("E7").Value = "  +3.38%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.174"
$ws.Range("E10").Value = "  +6.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.32"
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.494"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.14"
$ws.Range("E13").Value = "  +3.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000265"
$ws.Range("E14").Value = "  +4.68%  "
$ws.Range("D15").Value = "4.424.87"
$ws.Range("E15").Value = "  +3.25%  "
$ws.Range("D16").Value = "3.786.11"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("D17").Value = "70.287.45"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  +1.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.88"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "513.75"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.53"
$ws.Range("E22").Value = "  +3.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.729"
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.52"
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.14"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.22"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("E28").Value = "  +27.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.88"
$ws.Range("E31").Value = "  +4.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.76"
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.55"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.26"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("E37").Value = "  +4.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.135"
$ws.Range("E38").Value = "  +5.76%  "
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.15"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.05"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.02"
$ws.Range("E42").Value = "  -2.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.78"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "421.47"
$ws.Range("E44").Value = "  +5.92%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.85"
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "3.041.74"
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0366"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.57"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.79"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.47"
$ws.Range("E51").Value = "  +0.74%  "
